$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Remove the empty Title placeholder (Title 10) ---
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 13) {
        $shp.Delete()
    }
}

# --- Add the new credits/acknowledgement textbox ---
$box = $s.Shapes.AddTextbox(1, 57.1403937007874, 68.63039370078741, 621.5, 402.73913385826773)
$box.Name = "Text Placeholder 7"

$tr = $box.TextFrame.TextRange

$lines = @(
    'Please attribute Dr. Jim Alves-Foss and Dr. Jia Song, University of Idaho',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    'Except where otherwise noted, this work is licensed under https://creativecommons.org/licenses/by-nc-sa/4.0/',
    '',
    'Not withstanding the non-commercial license terms, non-profit educational institutions are granted a non-exclusive license to adapt and use this material, with attribution.',
    '',
    'Creative Commons and the double C in a circle are registered trademarks of Creative commons in the United States and other countries. Third party marks and brands are the property of their respective holders.',
    'Project sponsored by the National Security Agency under grant Number H98230-17-1-0199. The United States Government is authorized to reproduce and distribute reprints notwithstanding any copyright notation herein.',
    ''
)

$tr.Text = [string]::Join("`r", $lines)

# Paragraphs 9 through 15 (1-indexed) use 16pt type, matching the source deck.
$small = $tr.Paragraphs(9, 7)
$small.Font.Size = 16

Write-Output ("Shapes after edit: " + $s.Shapes.Count)
